$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Update the first three single-value rows ---
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"

# --- Insert 10 new single-value rows right after row 3 (before the old row 4) ---
$newValues = @("104","0.00002","0.00015","0.00008","0.00002","0.00009","0.00009","0.00015","0.00394","100.0")
$reversedValues = @($newValues[($newValues.Length-1)..0])

$beforeRow = $t.Rows(4)
foreach ($v in $reversedValues) {
    $newRow = $t.Rows.Add($beforeRow)
    $t.Cell($newRow.Index,1).Range.Text = $v
    $beforeRow = $newRow
}

# --- Collapse the three multi-value (tab-separated) rows down to a single value each ---
# After the 10-row insertion these are now rows 44, 45, 46 (were 34, 35, 36).
$t.Cell(44,1).Range.Text = "100"
$t.Cell(45,1).Range.Text = "0"
$t.Cell(46,1).Range.Text = "213"
